# Docker file and added new BookStore API
# - Adds a new "bookstoreUser" worksheet (after "PetStore") holding the
#   BookStore API test credentials (username/password).
# - Moves the selection/active cell on the "PetStore" sheet to J7 and makes
#   the new sheet the active tab.

$wb = $excel.ActiveWorkbook

# Record the new selection on the previously-active "PetStore" sheet before
# we move focus to the new sheet.
$petStore = $wb.Worksheets.Item(3)
$petStore.Range("J7").Select() | Out-Null

# Add the new sheet right after the last existing sheet ("PetStore").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$bookstoreUser = $wb.Worksheets.Add($null, $lastSheet)
$bookstoreUser.Name = "bookstoreUser"

# Header row.
$bookstoreUser.Range("A1").Value = "username"
$bookstoreUser.Range("B1").Value = "password"

# Data row - username kept as text (leading apostrophe -> quote-prefix
# style) and password cell turned into a hyperlink, matching the
# credentials-style row used on the other sheets (e.g. "PetStore"/G2).
$bookstoreUser.Range("A2").Value = "'ashuk"
$bookstoreUser.Range("B2").Value = "'Abcd@1234"
$bookstoreUser.Hyperlinks.Add($bookstoreUser.Range("B2"), "Abcd@1234") | Out-Null

# New sheet becomes the active tab, selection resting on A3.
$bookstoreUser.Range("A3").Select() | Out-Null
$bookstoreUser.Activate()
